$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C (CromosomaMax) to remain text so the long 0/1
# bit-strings are not coerced into numbers / scientific notation.
$ws.Range("C2:C201").NumberFormat = "@"

$ws.Range("C2").Value = "111011010011100111010101110010"
$ws.Range("D2").Value = 0.8587056948055914
$ws.Range("E2").Value = 0.00007779581916200958
$ws.Range("F2").Value = 0.2977361743313028
$ws.Range("C3").Value = "111011010011100111100111101100"
$ws.Range("D3").Value = 0.8587076728551368
$ws.Range("E3").Value = 0.01562082192970662
$ws.Range("F3").Value = 0.5982610274860772
$ws.Range("C4").Value = "111011010101001011100111101010"
$ws.Range("D4").Value = 0.8594148040659402
$ws.Range("E4").Value = 0.4650272860068329
$ws.Range("F4").Value = 0.8025423855049432
$ws.Range("C5").Value = "111011010101001011100111101010"
$ws.Range("D5").Value = 0.8594148040659402
$ws.Range("E5").Value = 0.8024748812251115
$ws.Range("F5").Value = 0.8365063847132603
$ws.Range("C6").Value = "111011010101001011100101110010"
$ws.Range("D6").Value = 0.8594145968549618
$ws.Range("E6").Value = 0.8024748845622623
$ws.Range("F6").Value = 0.8309587519584405
$ws.Range("C7").Value = "111011010101001011100111100100"
$ws.Range("D7").Value = 0.8594147937053906
$ws.Range("E7").Value = 0.8024746943446718
$ws.Range("F7").Value = 0.8364796941987744
$ws.Range("C8").Value = "111011010101001011100111101010"
$ws.Range("D8").Value = 0.8594148040659402
$ws.Range("E8").Value = 0.8024746943446718
$ws.Range("F8").Value = 0.8419661355251072
$ws.Range("C9").Value = "111011010101100111000111100010"
$ws.Range("D9").Value = 0.859609303305261
$ws.Range("E9").Value = 0.8024748812251115
$ws.Range("F9").Value = 0.8476060556579312
$ws.Range("C10").Value = "111011010011100111000111101010"
$ws.Range("D10").Value = 0.8587041344609606
$ws.Range("E10").Value = 0.8024748812251115
$ws.Range("F10").Value = 0.830819616617536
$ws.Range("C11").Value = "111011010011100111000111101010"
$ws.Range("D11").Value = 0.8587041344609606
$ws.Range("E11").Value = 0.8024748812251115
$ws.Range("F11").Value = 0.8419887712440488
$ws.Range("C12").Value = "111011010011100111000111101010"
$ws.Range("D12").Value = 0.8587041344609606
$ws.Range("E12").Value = 0.8024748812251115
$ws.Range("F12").Value = 0.8364425626536496
$ws.Range("C13").Value = "111011010101001111000111100010"
$ws.Range("D13").Value = 0.8594395452363079
$ws.Range("E13").Value = 0.8017642655163755
$ws.Range("F13").Value = 0.8308988645288986
$ws.Range("C14").Value = "111011010101001111000111100010"
$ws.Range("D14").Value = 0.8594395452363079
$ws.Range("E14").Value = 0.8017642655163755
$ws.Range("F14").Value = 0.841843805089362
$ws.Range("C15").Value = "111011010011111111000111101010"
$ws.Range("D15").Value = 0.858873819888117
$ws.Range("E15").Value = 0.8017642655163755
$ws.Range("F15").Value = 0.847464809823159
$ws.Range("C16").Value = "111011010011111111000111101010"
$ws.Range("D16").Value = 0.858873819888117
$ws.Range("E16").Value = 0.8587041344609606
$ws.Range("F16").Value = 0.8587212590381392
$ws.Range("C17").Value = "111011010011100111000111101010"
$ws.Range("D17").Value = 0.8587041344609606
$ws.Range("E17").Value = 0.8587041344609606
$ws.Range("F17").Value = 0.8587041344609606
$ws.Range("C18").Value = "111011010011100111000111101010"
$ws.Range("D18").Value = 0.8587041344609606
$ws.Range("E18").Value = 0.8587041344609606
$ws.Range("F18").Value = 0.8587041344609606
$ws.Range("C19").Value = "111011010011100111000111101010"
$ws.Range("D19").Value = 0.8587041344609606
$ws.Range("E19").Value = 0.8587041344609606
$ws.Range("F19").Value = 0.8587041344609606
$ws.Range("C20").Value = "111011010011100111000111101010"
$ws.Range("D20").Value = 0.8587041344609606
$ws.Range("E20").Value = 0.8587041344609606
$ws.Range("F20").Value = 0.8587041344609606
$ws.Range("C21").Value = "111011010011100111000111101010"
$ws.Range("D21").Value = 0.8587041344609606
$ws.Range("E21").Value = 0.8587041344609606
$ws.Range("F21").Value = 0.8587041344609606
$ws.Range("C22").Value = "111011010011100111000111101010"
$ws.Range("D22").Value = 0.8587041344609606
$ws.Range("E22").Value = 0.8587041344609606
$ws.Range("F22").Value = 0.8587041344609606
$ws.Range("C23").Value = "111011010011100111000111101010"
$ws.Range("D23").Value = 0.8587041344609606
$ws.Range("E23").Value = 0.8587041344609606
$ws.Range("F23").Value = 0.8587041344609606
$ws.Range("C24").Value = "111011010011100111000111101010"
$ws.Range("D24").Value = 0.8587041344609606
$ws.Range("E24").Value = 0.8587041344609606
$ws.Range("F24").Value = 0.8587041344609606
$ws.Range("C25").Value = "111011010011100111000111101010"
$ws.Range("D25").Value = 0.8587041344609606
$ws.Range("E25").Value = 0.8587041344609606
$ws.Range("F25").Value = 0.8587041344609606
$ws.Range("C26").Value = "111011010011100111000111101010"
$ws.Range("D26").Value = 0.8587041344609606
$ws.Range("E26").Value = 0.8587041344609606
$ws.Range("F26").Value = 0.8587041344609606
$ws.Range("C27").Value = "111011010011100111000111101010"
$ws.Range("D27").Value = 0.8587041344609606
$ws.Range("E27").Value = 0.8587041344609606
$ws.Range("F27").Value = 0.8587041344609606
$ws.Range("C28").Value = "111011010011100111000111101010"
$ws.Range("D28").Value = 0.8587041344609606
$ws.Range("E28").Value = 0.8587041344609606
$ws.Range("F28").Value = 0.8587041344609606
$ws.Range("C29").Value = "111011010011100111000111101010"
$ws.Range("D29").Value = 0.8587041344609606
$ws.Range("E29").Value = 0.8587041344609606
$ws.Range("F29").Value = 0.8587041344609606
$ws.Range("C30").Value = "111011010011100111000111101010"
$ws.Range("D30").Value = 0.8587041344609606
$ws.Range("E30").Value = 0.8587041344609606
$ws.Range("F30").Value = 0.8587041344609606
$ws.Range("C31").Value = "111011010011100111000111101010"
$ws.Range("D31").Value = 0.8587041344609606
$ws.Range("E31").Value = 0.8587041344609606
$ws.Range("F31").Value = 0.8587041344609606
$ws.Range("C32").Value = "111011010011100111000111101010"
$ws.Range("D32").Value = 0.8587041344609606
$ws.Range("E32").Value = 0.8587041344609606
$ws.Range("F32").Value = 0.8587041344609606
$ws.Range("C33").Value = "111011010011100111000111101010"
$ws.Range("D33").Value = 0.8587041344609606
$ws.Range("E33").Value = 0.8587041344609606
$ws.Range("F33").Value = 0.8587041344609606
$ws.Range("C34").Value = "111011010011100111000111101010"
$ws.Range("D34").Value = 0.8587041344609606
$ws.Range("E34").Value = 0.8587041344609606
$ws.Range("F34").Value = 0.8587041344609606
$ws.Range("C35").Value = "111011010011100111000111101010"
$ws.Range("D35").Value = 0.8587041344609606
$ws.Range("E35").Value = 0.8587041344609606
$ws.Range("F35").Value = 0.8587041344609606
$ws.Range("C36").Value = "111011010011100111000111101010"
$ws.Range("D36").Value = 0.8587041344609606
$ws.Range("E36").Value = 0.8587041344609606
$ws.Range("F36").Value = 0.8587041344609606
$ws.Range("C37").Value = "111011010011100111000111101010"
$ws.Range("D37").Value = 0.8587041344609606
$ws.Range("E37").Value = 0.8587041344609606
$ws.Range("F37").Value = 0.8587041344609606
$ws.Range("C38").Value = "111011010011100111000111101010"
$ws.Range("D38").Value = 0.8587041344609606
$ws.Range("E38").Value = 0.8587041344609606
$ws.Range("F38").Value = 0.8587041344609606
$ws.Range("C39").Value = "111011010011100111000111101010"
$ws.Range("D39").Value = 0.8587041344609606
$ws.Range("E39").Value = 0.8587041344609606
$ws.Range("F39").Value = 0.8587041344609606
$ws.Range("C40").Value = "111011010011100111000111101010"
$ws.Range("D40").Value = 0.8587041344609606
$ws.Range("E40").Value = 0.8587041344609606
$ws.Range("F40").Value = 0.8587041344609606
$ws.Range("C41").Value = "111011010011100111000111101010"
$ws.Range("D41").Value = 0.8587041344609606
$ws.Range("E41").Value = 0.8587041344609606
$ws.Range("F41").Value = 0.8587041344609606
$ws.Range("C42").Value = "111011010011100111000111101010"
$ws.Range("D42").Value = 0.8587041344609606
$ws.Range("E42").Value = 0.8587041344609606
$ws.Range("F42").Value = 0.8587041344609606
$ws.Range("C43").Value = "111011010011100111000111101010"
$ws.Range("D43").Value = 0.8587041344609606
$ws.Range("E43").Value = 0.8587041344609606
$ws.Range("F43").Value = 0.8587041344609606
$ws.Range("C44").Value = "111011010011100111000111101010"
$ws.Range("D44").Value = 0.8587041344609606
$ws.Range("E44").Value = 0.8587041344609606
$ws.Range("F44").Value = 0.8587041344609606
$ws.Range("C45").Value = "111011010011100111000111101010"
$ws.Range("D45").Value = 0.8587041344609606
$ws.Range("E45").Value = 0.8587041344609606
$ws.Range("F45").Value = 0.8587041344609606
$ws.Range("C46").Value = "111011010011100111000111101010"
$ws.Range("D46").Value = 0.8587041344609606
$ws.Range("E46").Value = 0.8587041344609606
$ws.Range("F46").Value = 0.8587041344609606
$ws.Range("C47").Value = "111011010011100111000111101010"
$ws.Range("D47").Value = 0.8587041344609606
$ws.Range("E47").Value = 0.8587041344609606
$ws.Range("F47").Value = 0.8587041344609606
$ws.Range("C48").Value = "111011010011100111000111101010"
$ws.Range("D48").Value = 0.8587041344609606
$ws.Range("E48").Value = 0.8587041344609606
$ws.Range("F48").Value = 0.8587041344609606
$ws.Range("C49").Value = "111011010011100111000111101010"
$ws.Range("D49").Value = 0.8587041344609606
$ws.Range("E49").Value = 0.8587041344609606
$ws.Range("F49").Value = 0.8587041344609606
$ws.Range("C50").Value = "111011010011100111000111101010"
$ws.Range("D50").Value = 0.8587041344609606
$ws.Range("E50").Value = 0.8587041344609606
$ws.Range("F50").Value = 0.8587041344609606
$ws.Range("C51").Value = "111011010011100111000111101010"
$ws.Range("D51").Value = 0.8587041344609606
$ws.Range("E51").Value = 0.8587041344609606
$ws.Range("F51").Value = 0.8587041344609606
$ws.Range("C52").Value = "111011010011100111000111101010"
$ws.Range("D52").Value = 0.8587041344609606
$ws.Range("E52").Value = 0.8587041344609606
$ws.Range("F52").Value = 0.8587041344609606
$ws.Range("C53").Value = "111011010011100111000111101010"
$ws.Range("D53").Value = 0.8587041344609606
$ws.Range("E53").Value = 0.8587041344609606
$ws.Range("F53").Value = 0.8587041344609606
$ws.Range("C54").Value = "111011010011100111000111101010"
$ws.Range("D54").Value = 0.8587041344609606
$ws.Range("E54").Value = 0.8587041344609606
$ws.Range("F54").Value = 0.8587041344609606
$ws.Range("C55").Value = "111011010011100111000111101010"
$ws.Range("D55").Value = 0.8587041344609606
$ws.Range("E55").Value = 0.8587041344609606
$ws.Range("F55").Value = 0.8587041344609606
$ws.Range("C56").Value = "111011010011100111000111101010"
$ws.Range("D56").Value = 0.8587041344609606
$ws.Range("E56").Value = 0.8587041344609606
$ws.Range("F56").Value = 0.8587041344609606
$ws.Range("C57").Value = "111011010011100111000111101010"
$ws.Range("D57").Value = 0.8587041344609606
$ws.Range("E57").Value = 0.8587041344609606
$ws.Range("F57").Value = 0.8587041344609606
$ws.Range("C58").Value = "111011010011100111000111101010"
$ws.Range("D58").Value = 0.8587041344609606
$ws.Range("E58").Value = 0.8587041344609606
$ws.Range("F58").Value = 0.8587041344609606
$ws.Range("C59").Value = "111011010011100111000111101010"
$ws.Range("D59").Value = 0.8587041344609606
$ws.Range("E59").Value = 0.8587041344609606
$ws.Range("F59").Value = 0.8587041344609606
$ws.Range("C60").Value = "111011010011100111000111101010"
$ws.Range("D60").Value = 0.8587041344609606
$ws.Range("E60").Value = 0.8587041344609606
$ws.Range("F60").Value = 0.8587041344609606
$ws.Range("C61").Value = "111011010011100111000111101010"
$ws.Range("D61").Value = 0.8587041344609606
$ws.Range("E61").Value = 0.8587041344609606
$ws.Range("F61").Value = 0.8587041344609606
$ws.Range("C62").Value = "111011010011100111000111101010"
$ws.Range("D62").Value = 0.8587041344609606
$ws.Range("E62").Value = 0.8587041344609606
$ws.Range("F62").Value = 0.8587041344609606
$ws.Range("C63").Value = "111011010011100111000111101010"
$ws.Range("D63").Value = 0.8587041344609606
$ws.Range("E63").Value = 0.8587041344609606
$ws.Range("F63").Value = 0.8587041344609606
$ws.Range("C64").Value = "111011010011100111000111101010"
$ws.Range("D64").Value = 0.8587041344609606
$ws.Range("E64").Value = 0.8587041344609606
$ws.Range("F64").Value = 0.8587041344609606
$ws.Range("C65").Value = "111011010011100111000111101010"
$ws.Range("D65").Value = 0.8587041344609606
$ws.Range("E65").Value = 0.8587041344609606
$ws.Range("F65").Value = 0.8587041344609606
$ws.Range("C66").Value = "111011010011100111000111101010"
$ws.Range("D66").Value = 0.8587041344609606
$ws.Range("E66").Value = 0.8587041344609606
$ws.Range("F66").Value = 0.8587041344609606
$ws.Range("C67").Value = "111011010011100111000111101010"
$ws.Range("D67").Value = 0.8587041344609606
$ws.Range("E67").Value = 0.8587041344609606
$ws.Range("F67").Value = 0.8587041344609606
$ws.Range("C68").Value = "111011010011100111000111101010"
$ws.Range("D68").Value = 0.8587041344609606
$ws.Range("E68").Value = 0.8587041344609606
$ws.Range("F68").Value = 0.8587041344609606
$ws.Range("C69").Value = "111011010011100111000111101010"
$ws.Range("D69").Value = 0.8587041344609606
$ws.Range("E69").Value = 0.8587041344609606
$ws.Range("F69").Value = 0.8587041344609606
$ws.Range("C70").Value = "111011010011100111000111101010"
$ws.Range("D70").Value = 0.8587041344609606
$ws.Range("E70").Value = 0.8587041344609606
$ws.Range("F70").Value = 0.8587041344609606
$ws.Range("C71").Value = "111011010011100111000111101010"
$ws.Range("D71").Value = 0.8587041344609606
$ws.Range("E71").Value = 0.8587041344609606
$ws.Range("F71").Value = 0.8587041344609606
$ws.Range("C72").Value = "111011010011100111000111101010"
$ws.Range("D72").Value = 0.8587041344609606
$ws.Range("E72").Value = 0.8587041344609606
$ws.Range("F72").Value = 0.8587041344609606
$ws.Range("C73").Value = "111011010011100111000111101010"
$ws.Range("D73").Value = 0.8587041344609606
$ws.Range("E73").Value = 0.8587041344609606
$ws.Range("F73").Value = 0.8587041344609606
$ws.Range("C74").Value = "111011010011100111000111101010"
$ws.Range("D74").Value = 0.8587041344609606
$ws.Range("E74").Value = 0.8587041344609606
$ws.Range("F74").Value = 0.8587041344609606
$ws.Range("C75").Value = "111011010011100111000111101010"
$ws.Range("D75").Value = 0.8587041344609606
$ws.Range("E75").Value = 0.8587041344609606
$ws.Range("F75").Value = 0.8587041344609606
$ws.Range("C76").Value = "111011010011100111000111101010"
$ws.Range("D76").Value = 0.8587041344609606
$ws.Range("E76").Value = 0.8587041344609606
$ws.Range("F76").Value = 0.8587041344609606
$ws.Range("C77").Value = "111011010011100111000111101010"
$ws.Range("D77").Value = 0.8587041344609606
$ws.Range("E77").Value = 0.8587041344609606
$ws.Range("F77").Value = 0.8587041344609606
$ws.Range("C78").Value = "111011010011100111000111101010"
$ws.Range("D78").Value = 0.8587041344609606
$ws.Range("E78").Value = 0.8587041344609606
$ws.Range("F78").Value = 0.8587041344609606
$ws.Range("C79").Value = "111011010011100111000111101010"
$ws.Range("D79").Value = 0.8587041344609606
$ws.Range("E79").Value = 0.8587041344609606
$ws.Range("F79").Value = 0.8587041344609606
$ws.Range("C80").Value = "111011010011100111000111101010"
$ws.Range("D80").Value = 0.8587041344609606
$ws.Range("E80").Value = 0.8587041344609606
$ws.Range("F80").Value = 0.8587041344609606
$ws.Range("C81").Value = "111011010011100111000111101010"
$ws.Range("D81").Value = 0.8587041344609606
$ws.Range("E81").Value = 0.8587041344609606
$ws.Range("F81").Value = 0.8587041344609606
$ws.Range("C82").Value = "111011010011100111000111101010"
$ws.Range("D82").Value = 0.8587041344609606
$ws.Range("E82").Value = 0.8587041344609606
$ws.Range("F82").Value = 0.8587041344609606
$ws.Range("C83").Value = "111011010011100111000111101010"
$ws.Range("D83").Value = 0.8587041344609606
$ws.Range("E83").Value = 0.8587041344609606
$ws.Range("F83").Value = 0.8587041344609606
$ws.Range("C84").Value = "111011010011100111000111101010"
$ws.Range("D84").Value = 0.8587041344609606
$ws.Range("E84").Value = 0.8587041344609606
$ws.Range("F84").Value = 0.8587041344609606
$ws.Range("C85").Value = "111011010011100111000111101010"
$ws.Range("D85").Value = 0.8587041344609606
$ws.Range("E85").Value = 0.8587041344609606
$ws.Range("F85").Value = 0.8587041344609606
$ws.Range("C86").Value = "111011010011100111000111101010"
$ws.Range("D86").Value = 0.8587041344609606
$ws.Range("E86").Value = 0.8587041344609606
$ws.Range("F86").Value = 0.8587041344609606
$ws.Range("C87").Value = "111011010011100111000111101010"
$ws.Range("D87").Value = 0.8587041344609606
$ws.Range("E87").Value = 0.8587041344609606
$ws.Range("F87").Value = 0.8587041344609606
$ws.Range("C88").Value = "111011010011100111000111101010"
$ws.Range("D88").Value = 0.8587041344609606
$ws.Range("E88").Value = 0.8587041344609606
$ws.Range("F88").Value = 0.8587041344609606
$ws.Range("C89").Value = "111011010011100111000111101010"
$ws.Range("D89").Value = 0.8587041344609606
$ws.Range("E89").Value = 0.8587041344609606
$ws.Range("F89").Value = 0.8587041344609606
$ws.Range("C90").Value = "111011010011100111000111101010"
$ws.Range("D90").Value = 0.8587041344609606
$ws.Range("E90").Value = 0.8587041344609606
$ws.Range("F90").Value = 0.8587041344609606
$ws.Range("C91").Value = "111011010011100111000111101010"
$ws.Range("D91").Value = 0.8587041344609606
$ws.Range("E91").Value = 0.8587041344609606
$ws.Range("F91").Value = 0.8587041344609606
$ws.Range("C92").Value = "111011010011100111000111101010"
$ws.Range("D92").Value = 0.8587041344609606
$ws.Range("E92").Value = 0.8587041344609606
$ws.Range("F92").Value = 0.8587041344609606
$ws.Range("C93").Value = "111011010011100111000111101010"
$ws.Range("D93").Value = 0.8587041344609606
$ws.Range("E93").Value = 0.8587041344609606
$ws.Range("F93").Value = 0.8587041344609606
$ws.Range("C94").Value = "111011010011100111000111101010"
$ws.Range("D94").Value = 0.8587041344609606
$ws.Range("E94").Value = 0.8587041344609606
$ws.Range("F94").Value = 0.8587041344609606
$ws.Range("C95").Value = "111011010011100111000111101010"
$ws.Range("D95").Value = 0.8587041344609606
$ws.Range("E95").Value = 0.8587041344609606
$ws.Range("F95").Value = 0.8587041344609606
$ws.Range("C96").Value = "111011010011100111000111101010"
$ws.Range("D96").Value = 0.8587041344609606
$ws.Range("E96").Value = 0.8587041344609606
$ws.Range("F96").Value = 0.8587041344609606
$ws.Range("C97").Value = "111011010011100111000111101010"
$ws.Range("D97").Value = 0.8587041344609606
$ws.Range("E97").Value = 0.8587041344609606
$ws.Range("F97").Value = 0.8587041344609606
$ws.Range("C98").Value = "111011010011100111000111101010"
$ws.Range("D98").Value = 0.8587041344609606
$ws.Range("E98").Value = 0.8587041344609606
$ws.Range("F98").Value = 0.8587041344609606
$ws.Range("C99").Value = "111011010011100111000111101010"
$ws.Range("D99").Value = 0.8587041344609606
$ws.Range("E99").Value = 0.8587041344609606
$ws.Range("F99").Value = 0.8587041344609606
$ws.Range("C100").Value = "111011010011100111000111101010"
$ws.Range("D100").Value = 0.8587041344609606
$ws.Range("E100").Value = 0.8587041344609606
$ws.Range("F100").Value = 0.8587041344609606
$ws.Range("C101").Value = "111011010011100111000111101010"
$ws.Range("D101").Value = 0.8587041344609606
$ws.Range("E101").Value = 0.8587041344609606
$ws.Range("F101").Value = 0.8587041344609606
$ws.Range("C102").Value = "111011010011100111000111101010"
$ws.Range("D102").Value = 0.8587041344609606
$ws.Range("E102").Value = 0.8587041344609606
$ws.Range("F102").Value = 0.8587041344609606
$ws.Range("C103").Value = "111011010011100111000111101010"
$ws.Range("D103").Value = 0.8587041344609606
$ws.Range("E103").Value = 0.8587041344609606
$ws.Range("F103").Value = 0.8587041344609606
$ws.Range("C104").Value = "111011010011100111000111101010"
$ws.Range("D104").Value = 0.8587041344609606
$ws.Range("E104").Value = 0.8587041344609606
$ws.Range("F104").Value = 0.8587041344609606
$ws.Range("C105").Value = "111011010011100111000111101010"
$ws.Range("D105").Value = 0.8587041344609606
$ws.Range("E105").Value = 0.8587041344609606
$ws.Range("F105").Value = 0.8587041344609606
$ws.Range("C106").Value = "111011010011100111000111101010"
$ws.Range("D106").Value = 0.8587041344609606
$ws.Range("E106").Value = 0.8587041344609606
$ws.Range("F106").Value = 0.8587041344609606
$ws.Range("C107").Value = "111011010011100111000111101010"
$ws.Range("D107").Value = 0.8587041344609606
$ws.Range("E107").Value = 0.8587041344609606
$ws.Range("F107").Value = 0.8587041344609606
$ws.Range("C108").Value = "111011010011100111000111101010"
$ws.Range("D108").Value = 0.8587041344609606
$ws.Range("E108").Value = 0.8587041344609606
$ws.Range("F108").Value = 0.8587041344609606
$ws.Range("C109").Value = "111011010011100111000111101010"
$ws.Range("D109").Value = 0.8587041344609606
$ws.Range("E109").Value = 0.8587041344609606
$ws.Range("F109").Value = 0.8587041344609606
$ws.Range("C110").Value = "111011010011100111000111101010"
$ws.Range("D110").Value = 0.8587041344609606
$ws.Range("E110").Value = 0.8587041344609606
$ws.Range("F110").Value = 0.8587041344609606
$ws.Range("C111").Value = "111011010011100111000111101010"
$ws.Range("D111").Value = 0.8587041344609606
$ws.Range("E111").Value = 0.8587041344609606
$ws.Range("F111").Value = 0.8587041344609606
$ws.Range("C112").Value = "111011010011100111000111101010"
$ws.Range("D112").Value = 0.8587041344609606
$ws.Range("E112").Value = 0.8587041344609606
$ws.Range("F112").Value = 0.8587041344609606
$ws.Range("C113").Value = "111011010011100111000111101010"
$ws.Range("D113").Value = 0.8587041344609606
$ws.Range("E113").Value = 0.8587041344609606
$ws.Range("F113").Value = 0.8587041344609606
$ws.Range("C114").Value = "111011010011100111000111101010"
$ws.Range("D114").Value = 0.8587041344609606
$ws.Range("E114").Value = 0.8587041344609606
$ws.Range("F114").Value = 0.8587041344609606
$ws.Range("C115").Value = "111011010011100111000111101010"
$ws.Range("D115").Value = 0.8587041344609606
$ws.Range("E115").Value = 0.8587041344609606
$ws.Range("F115").Value = 0.8587041344609606
$ws.Range("C116").Value = "111011010011100111000111101010"
$ws.Range("D116").Value = 0.8587041344609606
$ws.Range("E116").Value = 0.8587041344609606
$ws.Range("F116").Value = 0.8587041344609606
$ws.Range("C117").Value = "111011010011100111000111101010"
$ws.Range("D117").Value = 0.8587041344609606
$ws.Range("E117").Value = 0.8587041344609606
$ws.Range("F117").Value = 0.8587041344609606
$ws.Range("C118").Value = "111011010011100111000111101010"
$ws.Range("D118").Value = 0.8587041344609606
$ws.Range("E118").Value = 0.8587041344609606
$ws.Range("F118").Value = 0.8587041344609606
$ws.Range("C119").Value = "111011010011100111000111101010"
$ws.Range("D119").Value = 0.8587041344609606
$ws.Range("E119").Value = 0.8587041344609606
$ws.Range("F119").Value = 0.8587041344609606
$ws.Range("C120").Value = "111011010011100111000111101010"
$ws.Range("D120").Value = 0.8587041344609606
$ws.Range("E120").Value = 0.8587041344609606
$ws.Range("F120").Value = 0.8587041344609606
$ws.Range("C121").Value = "111011010011100111000111101010"
$ws.Range("D121").Value = 0.8587041344609606
$ws.Range("E121").Value = 0.8587041344609606
$ws.Range("F121").Value = 0.8587041344609606
$ws.Range("C122").Value = "111011010011100111000111101010"
$ws.Range("D122").Value = 0.8587041344609606
$ws.Range("E122").Value = 0.8587041344609606
$ws.Range("F122").Value = 0.8587041344609606
$ws.Range("C123").Value = "111011010011100111000111101010"
$ws.Range("D123").Value = 0.8587041344609606
$ws.Range("E123").Value = 0.8587041344609606
$ws.Range("F123").Value = 0.8587041344609606
$ws.Range("C124").Value = "111011010011100111000111101010"
$ws.Range("D124").Value = 0.8587041344609606
$ws.Range("E124").Value = 0.8587041344609606
$ws.Range("F124").Value = 0.8587041344609606
$ws.Range("C125").Value = "111011010011100111000111101010"
$ws.Range("D125").Value = 0.8587041344609606
$ws.Range("E125").Value = 0.8587041344609606
$ws.Range("F125").Value = 0.8587041344609606
$ws.Range("C126").Value = "111011010011100111000111101010"
$ws.Range("D126").Value = 0.8587041344609606
$ws.Range("E126").Value = 0.8587041344609606
$ws.Range("F126").Value = 0.8587041344609606
$ws.Range("C127").Value = "111011010011100111000111101010"
$ws.Range("D127").Value = 0.8587041344609606
$ws.Range("E127").Value = 0.8587041344609606
$ws.Range("F127").Value = 0.8587041344609606
$ws.Range("C128").Value = "111011010011100111000111101010"
$ws.Range("D128").Value = 0.8587041344609606
$ws.Range("E128").Value = 0.8587041344609606
$ws.Range("F128").Value = 0.8587041344609606
$ws.Range("C129").Value = "111011010011100111000111101010"
$ws.Range("D129").Value = 0.8587041344609606
$ws.Range("E129").Value = 0.8587041344609606
$ws.Range("F129").Value = 0.8587041344609606
$ws.Range("C130").Value = "111011010011100111000111101010"
$ws.Range("D130").Value = 0.8587041344609606
$ws.Range("E130").Value = 0.8587041344609606
$ws.Range("F130").Value = 0.8587041344609606
$ws.Range("C131").Value = "111011010011100111000111101010"
$ws.Range("D131").Value = 0.8587041344609606
$ws.Range("E131").Value = 0.8587041344609606
$ws.Range("F131").Value = 0.8587041344609606
$ws.Range("C132").Value = "111011010011100111000111101010"
$ws.Range("D132").Value = 0.8587041344609606
$ws.Range("E132").Value = 0.8587041344609606
$ws.Range("F132").Value = 0.8587041344609606
$ws.Range("C133").Value = "111011010011100111000111101010"
$ws.Range("D133").Value = 0.8587041344609606
$ws.Range("E133").Value = 0.8587041344609606
$ws.Range("F133").Value = 0.8587041344609606
$ws.Range("C134").Value = "111011010011100111000111101010"
$ws.Range("D134").Value = 0.8587041344609606
$ws.Range("E134").Value = 0.8587041344609606
$ws.Range("F134").Value = 0.8587041344609606
$ws.Range("C135").Value = "111011010011100111000111101010"
$ws.Range("D135").Value = 0.8587041344609606
$ws.Range("E135").Value = 0.8587041344609606
$ws.Range("F135").Value = 0.8587041344609606
$ws.Range("C136").Value = "111011010011100111000111101010"
$ws.Range("D136").Value = 0.8587041344609606
$ws.Range("E136").Value = 0.8587041344609606
$ws.Range("F136").Value = 0.8587041344609606
$ws.Range("C137").Value = "111011010011100111000111101010"
$ws.Range("D137").Value = 0.8587041344609606
$ws.Range("E137").Value = 0.8587041344609606
$ws.Range("F137").Value = 0.8587041344609606
$ws.Range("C138").Value = "111011010011100111000111101010"
$ws.Range("D138").Value = 0.8587041344609606
$ws.Range("E138").Value = 0.8587041344609606
$ws.Range("F138").Value = 0.8587041344609606
$ws.Range("C139").Value = "111011010011100111000111101010"
$ws.Range("D139").Value = 0.8587041344609606
$ws.Range("E139").Value = 0.8587041344609606
$ws.Range("F139").Value = 0.8587041344609606
$ws.Range("C140").Value = "111011010011100111000111101010"
$ws.Range("D140").Value = 0.8587041344609606
$ws.Range("E140").Value = 0.8587041344609606
$ws.Range("F140").Value = 0.8587041344609606
$ws.Range("C141").Value = "111011010011100111000111101010"
$ws.Range("D141").Value = 0.8587041344609606
$ws.Range("E141").Value = 0.8587041344609606
$ws.Range("F141").Value = 0.8587041344609606
$ws.Range("C142").Value = "111011010011100111000111101010"
$ws.Range("D142").Value = 0.8587041344609606
$ws.Range("E142").Value = 0.8587041344609606
$ws.Range("F142").Value = 0.8587041344609606
$ws.Range("C143").Value = "111011010011100111000111101010"
$ws.Range("D143").Value = 0.8587041344609606
$ws.Range("E143").Value = 0.8587041344609606
$ws.Range("F143").Value = 0.8587041344609606
$ws.Range("C144").Value = "111011010011100111000111101010"
$ws.Range("D144").Value = 0.8587041344609606
$ws.Range("E144").Value = 0.8587041344609606
$ws.Range("F144").Value = 0.8587041344609606
$ws.Range("C145").Value = "111011010011100111000111101010"
$ws.Range("D145").Value = 0.8587041344609606
$ws.Range("E145").Value = 0.8587041344609606
$ws.Range("F145").Value = 0.8587041344609606
$ws.Range("C146").Value = "111011010011100111000111101010"
$ws.Range("D146").Value = 0.8587041344609606
$ws.Range("E146").Value = 0.8587041344609606
$ws.Range("F146").Value = 0.8587041344609606
$ws.Range("C147").Value = "111011010011100111000111101010"
$ws.Range("D147").Value = 0.8587041344609606
$ws.Range("E147").Value = 0.8587041344609606
$ws.Range("F147").Value = 0.8587041344609606
$ws.Range("C148").Value = "111011010011100111000111101010"
$ws.Range("D148").Value = 0.8587041344609606
$ws.Range("E148").Value = 0.8587041344609606
$ws.Range("F148").Value = 0.8587041344609606
$ws.Range("C149").Value = "111011010011100111000111101010"
$ws.Range("D149").Value = 0.8587041344609606
$ws.Range("E149").Value = 0.8587041344609606
$ws.Range("F149").Value = 0.8587041344609606
$ws.Range("C150").Value = "111011010011100111000111101010"
$ws.Range("D150").Value = 0.8587041344609606
$ws.Range("E150").Value = 0.8587041344609606
$ws.Range("F150").Value = 0.8587041344609606
$ws.Range("C151").Value = "111011010011100111000111101010"
$ws.Range("D151").Value = 0.8587041344609606
$ws.Range("E151").Value = 0.8587041344609606
$ws.Range("F151").Value = 0.8587041344609606
$ws.Range("C152").Value = "111011010011100111000111101010"
$ws.Range("D152").Value = 0.8587041344609606
$ws.Range("E152").Value = 0.8587041344609606
$ws.Range("F152").Value = 0.8587041344609606
$ws.Range("C153").Value = "111011010011100111000111101010"
$ws.Range("D153").Value = 0.8587041344609606
$ws.Range("E153").Value = 0.8587041344609606
$ws.Range("F153").Value = 0.8587041344609606
$ws.Range("C154").Value = "111011010011100111000111101010"
$ws.Range("D154").Value = 0.8587041344609606
$ws.Range("E154").Value = 0.8587041344609606
$ws.Range("F154").Value = 0.8587041344609606
$ws.Range("C155").Value = "111011010011100111000111101010"
$ws.Range("D155").Value = 0.8587041344609606
$ws.Range("E155").Value = 0.8587041344609606
$ws.Range("F155").Value = 0.8587041344609606
$ws.Range("C156").Value = "111011010011100111000111101010"
$ws.Range("D156").Value = 0.8587041344609606
$ws.Range("E156").Value = 0.8587041344609606
$ws.Range("F156").Value = 0.8587041344609606
$ws.Range("C157").Value = "111011010011100111000111101010"
$ws.Range("D157").Value = 0.8587041344609606
$ws.Range("E157").Value = 0.8587041344609606
$ws.Range("F157").Value = 0.8587041344609606
$ws.Range("C158").Value = "111011010011100111000111101010"
$ws.Range("D158").Value = 0.8587041344609606
$ws.Range("E158").Value = 0.8587041344609606
$ws.Range("F158").Value = 0.8587041344609606
$ws.Range("C159").Value = "111011010011100111000111101010"
$ws.Range("D159").Value = 0.8587041344609606
$ws.Range("E159").Value = 0.8587041344609606
$ws.Range("F159").Value = 0.8587041344609606
$ws.Range("C160").Value = "111011010011100111000111101010"
$ws.Range("D160").Value = 0.8587041344609606
$ws.Range("E160").Value = 0.8587041344609606
$ws.Range("F160").Value = 0.8587041344609606
$ws.Range("C161").Value = "111011010011100111000111101010"
$ws.Range("D161").Value = 0.8587041344609606
$ws.Range("E161").Value = 0.8587041344609606
$ws.Range("F161").Value = 0.8587041344609606
$ws.Range("C162").Value = "111011010011100111000111101010"
$ws.Range("D162").Value = 0.8587041344609606
$ws.Range("E162").Value = 0.8587041344609606
$ws.Range("F162").Value = 0.8587041344609606
$ws.Range("C163").Value = "111011010011100111000111101010"
$ws.Range("D163").Value = 0.8587041344609606
$ws.Range("E163").Value = 0.8587041344609606
$ws.Range("F163").Value = 0.8587041344609606
$ws.Range("C164").Value = "111011010011100111000111101010"
$ws.Range("D164").Value = 0.8587041344609606
$ws.Range("E164").Value = 0.8587041344609606
$ws.Range("F164").Value = 0.8587041344609606
$ws.Range("C165").Value = "111011010011100111000111101010"
$ws.Range("D165").Value = 0.8587041344609606
$ws.Range("E165").Value = 0.8587041344609606
$ws.Range("F165").Value = 0.8587041344609606
$ws.Range("C166").Value = "111011010011100111000111101010"
$ws.Range("D166").Value = 0.8587041344609606
$ws.Range("E166").Value = 0.8587041344609606
$ws.Range("F166").Value = 0.8587041344609606
$ws.Range("C167").Value = "111011010011100111000111101010"
$ws.Range("D167").Value = 0.8587041344609606
$ws.Range("E167").Value = 0.8587041344609606
$ws.Range("F167").Value = 0.8587041344609606
$ws.Range("C168").Value = "111011010011100111000111101010"
$ws.Range("D168").Value = 0.8587041344609606
$ws.Range("E168").Value = 0.8587041344609606
$ws.Range("F168").Value = 0.8587041344609606
$ws.Range("C169").Value = "111011010011100111000111101010"
$ws.Range("D169").Value = 0.8587041344609606
$ws.Range("E169").Value = 0.8587041344609606
$ws.Range("F169").Value = 0.8587041344609606
$ws.Range("C170").Value = "111011010011100111000111101010"
$ws.Range("D170").Value = 0.8587041344609606
$ws.Range("E170").Value = 0.8587041344609606
$ws.Range("F170").Value = 0.8587041344609606
$ws.Range("C171").Value = "111011010011100111000111101010"
$ws.Range("D171").Value = 0.8587041344609606
$ws.Range("E171").Value = 0.8587041344609606
$ws.Range("F171").Value = 0.8587041344609606
$ws.Range("C172").Value = "111011010011100111000111101010"
$ws.Range("D172").Value = 0.8587041344609606
$ws.Range("E172").Value = 0.8587041344609606
$ws.Range("F172").Value = 0.8587041344609606
$ws.Range("C173").Value = "111011010011100111000111101010"
$ws.Range("D173").Value = 0.8587041344609606
$ws.Range("E173").Value = 0.8587041344609606
$ws.Range("F173").Value = 0.8587041344609606
$ws.Range("C174").Value = "111011010011100111000111101010"
$ws.Range("D174").Value = 0.8587041344609606
$ws.Range("E174").Value = 0.8587041344609606
$ws.Range("F174").Value = 0.8587041344609606
$ws.Range("C175").Value = "111011010011100111000111101010"
$ws.Range("D175").Value = 0.8587041344609606
$ws.Range("E175").Value = 0.8587041344609606
$ws.Range("F175").Value = 0.8587041344609606
$ws.Range("C176").Value = "111011010011100111000111101010"
$ws.Range("D176").Value = 0.8587041344609606
$ws.Range("E176").Value = 0.8587041344609606
$ws.Range("F176").Value = 0.8587041344609606
$ws.Range("C177").Value = "111011010011100111000111101010"
$ws.Range("D177").Value = 0.8587041344609606
$ws.Range("E177").Value = 0.8587041344609606
$ws.Range("F177").Value = 0.8587041344609606
$ws.Range("C178").Value = "111011010011100111000111101010"
$ws.Range("D178").Value = 0.8587041344609606
$ws.Range("E178").Value = 0.8587041344609606
$ws.Range("F178").Value = 0.8587041344609606
$ws.Range("C179").Value = "111011010011100111000111101010"
$ws.Range("D179").Value = 0.8587041344609606
$ws.Range("E179").Value = 0.8587041344609606
$ws.Range("F179").Value = 0.8587041344609606
$ws.Range("C180").Value = "111011010011100111000111101010"
$ws.Range("D180").Value = 0.8587041344609606
$ws.Range("E180").Value = 0.8587041344609606
$ws.Range("F180").Value = 0.8587041344609606
$ws.Range("C181").Value = "111011010011100111000111101010"
$ws.Range("D181").Value = 0.8587041344609606
$ws.Range("E181").Value = 0.8587041344609606
$ws.Range("F181").Value = 0.8587041344609606
$ws.Range("C182").Value = "111011010011100111000111101010"
$ws.Range("D182").Value = 0.8587041344609606
$ws.Range("E182").Value = 0.8587041344609606
$ws.Range("F182").Value = 0.8587041344609606
$ws.Range("C183").Value = "111011010011100111000111101010"
$ws.Range("D183").Value = 0.8587041344609606
$ws.Range("E183").Value = 0.8587041344609606
$ws.Range("F183").Value = 0.8587041344609606
$ws.Range("C184").Value = "111011010011100111000111101010"
$ws.Range("D184").Value = 0.8587041344609606
$ws.Range("E184").Value = 0.8587041344609606
$ws.Range("F184").Value = 0.8587041344609606
$ws.Range("C185").Value = "111011010011100111000111101010"
$ws.Range("D185").Value = 0.8587041344609606
$ws.Range("E185").Value = 0.8587041344609606
$ws.Range("F185").Value = 0.8587041344609606
$ws.Range("C186").Value = "111011010011100111000111101010"
$ws.Range("D186").Value = 0.8587041344609606
$ws.Range("E186").Value = 0.8587041344609606
$ws.Range("F186").Value = 0.8587041344609606
$ws.Range("C187").Value = "111011010011100111000111101010"
$ws.Range("D187").Value = 0.8587041344609606
$ws.Range("E187").Value = 0.8587041344609606
$ws.Range("F187").Value = 0.8587041344609606
$ws.Range("C188").Value = "111011010011100111000111101010"
$ws.Range("D188").Value = 0.8587041344609606
$ws.Range("E188").Value = 0.8587041344609606
$ws.Range("F188").Value = 0.8587041344609606
$ws.Range("C189").Value = "111011010011100111000111101010"
$ws.Range("D189").Value = 0.8587041344609606
$ws.Range("E189").Value = 0.8587041344609606
$ws.Range("F189").Value = 0.8587041344609606
$ws.Range("C190").Value = "111011010011100111000111101010"
$ws.Range("D190").Value = 0.8587041344609606
$ws.Range("E190").Value = 0.8587041344609606
$ws.Range("F190").Value = 0.8587041344609606
$ws.Range("C191").Value = "111011010011100111000111101010"
$ws.Range("D191").Value = 0.8587041344609606
$ws.Range("E191").Value = 0.8587041344609606
$ws.Range("F191").Value = 0.8587041344609606
$ws.Range("C192").Value = "111011010011100111000111101010"
$ws.Range("D192").Value = 0.8587041344609606
$ws.Range("E192").Value = 0.8587041344609606
$ws.Range("F192").Value = 0.8587041344609606
$ws.Range("C193").Value = "111011010011100111000111101010"
$ws.Range("D193").Value = 0.8587041344609606
$ws.Range("E193").Value = 0.8587041344609606
$ws.Range("F193").Value = 0.8587041344609606
$ws.Range("C194").Value = "111011010011100111000111101010"
$ws.Range("D194").Value = 0.8587041344609606
$ws.Range("E194").Value = 0.8587041344609606
$ws.Range("F194").Value = 0.8587041344609606
$ws.Range("C195").Value = "111011010011100111000111101010"
$ws.Range("D195").Value = 0.8587041344609606
$ws.Range("E195").Value = 0.8587041344609606
$ws.Range("F195").Value = 0.8587041344609606
$ws.Range("C196").Value = "111011010011100111000111101010"
$ws.Range("D196").Value = 0.8587041344609606
$ws.Range("E196").Value = 0.8587041344609606
$ws.Range("F196").Value = 0.8587041344609606
$ws.Range("C197").Value = "111011010011100111000111101010"
$ws.Range("D197").Value = 0.8587041344609606
$ws.Range("E197").Value = 0.8587041344609606
$ws.Range("F197").Value = 0.8587041344609606
$ws.Range("C198").Value = "111011010011100111000111101010"
$ws.Range("D198").Value = 0.8587041344609606
$ws.Range("E198").Value = 0.8587041344609606
$ws.Range("F198").Value = 0.8587041344609606
$ws.Range("C199").Value = "111011010011100111000111101010"
$ws.Range("D199").Value = 0.8587041344609606
$ws.Range("E199").Value = 0.8587041344609606
$ws.Range("F199").Value = 0.8587041344609606
$ws.Range("C200").Value = "111011010011100111000111101010"
$ws.Range("D200").Value = 0.8587041344609606
$ws.Range("E200").Value = 0.8587041344609606
$ws.Range("F200").Value = 0.8587041344609606
$ws.Range("C201").Value = "111011010011100111000111101010"
$ws.Range("D201").Value = 0.8587041344609606
$ws.Range("E201").Value = 0.8587041344609606
$ws.Range("F201").Value = 0.8587041344609606
